# Case with 380 kV: refresh the vm_pu results for the slack bus voltage scenario
# (B column changed from 1.05 -> 1.02 p.u.; all dependent bus voltages recomputed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, in column order B,C,D,E,F and I,J,K,L,M,N (G/H are untouched).
$rowsData = @(
    @(1.02, 1.025592242350629, 1.029703073679386, 0.9926147277508489, 1.024094311281119, 1.032349696484665, 1.030760403123009, 1.03251621596089, 0.9955398523336033, 1.026923821697653, 1.032224201125319),
    @(1.02, 1.026428456599205, 1.030303031458803, 0.9936372048519304, 1.025564042921461, 1.032513588753615, 1.031236601235201, 1.032924935151167, 0.9963617723202692, 1.028198768326416, 1.032701075493446),
    @(1.02, 1.026969553194661, 1.03069125267988, 0.9942998659930995, 1.026515113353726, 1.032618460751116, 1.031544099509201, 1.033188717207306, 0.9968940712668345, 1.029023303398076, 1.033009010450266),
    @(1.02, 1.027197032078154, 1.030854461960202, 0.9945786998346017, 1.026914960621605, 1.032662266843601, 1.031673219758633, 1.033299446350823, 0.997117960005301, 1.029369835723623, 1.033138314065262),
    @(1.02, 1.027235226878024, 1.030881865552444, 0.9946255319796338, 1.026982097862731, 1.032669605529604, 1.031694890692937, 1.033318028581882, 0.9971555583673453, 1.029428014126861, 1.033160015774779),
    @(1.02, 1.026972592771941, 1.03069343348722, 0.9943035907982488, 1.026520456059809, 1.032619047199031, 1.031545825418036, 1.033190197423808, 0.9968970624462087, 1.029027934177543, 1.033010738810089),
    @(1.02, 1.025874842468228, 1.02990582991538, 0.9929600610674301, 1.02459100477046, 1.032405328249536, 1.030921467544184, 1.032654486411489, 0.995817528259106, 1.027354789091173, 1.032385494276441),
    @(1.02, 1.023940566282738, 1.028518077207004, 0.9906006454969559, 1.021191317514622, 1.032019722924527, 1.029816428456417, 1.031705256719966, 0.9939188001724441, 1.024402972133717, 1.031278885906421),
    @(1.02, 1.022651144897556, 1.027593033464755, 0.989033133672735, 1.018924783183762, 1.031756612378359, 1.029076496880997, 1.031068940524896, 0.9926553831429383, 1.022432519768097, 1.030537903543292),
    @(1.02, 1.022092839306676, 1.027192519834161, 0.988355674866747, 1.017943275545737, 1.031641252663634, 1.02875533289007, 1.030792584249455, 0.9921088820399291, 1.021578635017787, 1.03021628346266),
    @(1.02, 1.021885463454856, 1.027043757546277, 0.9881042295826724, 1.017578682562724, 1.031598188097417, 1.028635922941443, 1.030689809294589, 0.9919059725120875, 1.021261360075911, 1.030096703938209),
    @(1.02, 1.02192994610602, 1.027075667263243, 0.9881581567098651, 1.017656889857775, 1.031607435311825, 1.028661541997886, 1.030711860473901, 0.9919494934313052, 1.02132942140565, 1.030122359376652),
    @(1.02, 1.022075697476571, 1.027180222963559, 0.9883348863814464, 1.017913138557826, 1.031637697314641, 1.028745464778561, 1.030784091365782, 0.9920921077337197, 1.021552411119946, 1.030206401337301),
    @(1.02, 1.022165500207211, 1.027244644041976, 0.9884438009545853, 1.018071019396213, 1.031656314271319, 1.028797157093973, 1.030828578790511, 0.9921799884222134, 1.021689788522296, 1.03025816706173),
    @(1.02, 1.022688198301587, 1.027619615074151, 0.9890781214508737, 1.01898992044611, 1.031764238285915, 1.029097795293988, 1.031087263982699, 0.9926916645766087, 1.022489174909079, 1.030559232202473),
    @(1.02, 1.023016079240298, 1.027854834838296, 0.989476357848556, 1.019566297089691, 1.031831553144755, 1.029286171990116, 1.031249309157577, 0.9930127773699352, 1.022990427717103, 1.030747876415117),
    @(1.02, 1.023207328781069, 1.027992038004897, 0.9897087662937556, 1.01990248012018, 1.031870678669711, 1.029395974748137, 1.03134374753446, 0.9932001317071769, 1.023282736027173, 1.030857835105648),
    @(1.02, 1.023272540284078, 1.028038821299989, 0.9897880325774034, 1.020017108683994, 1.03188399603421, 1.029433402063876, 1.031375935029077, 0.9932640239640975, 1.023382394880294, 1.030895315572468),
    @(1.02, 1.02298090049505, 1.027829597619938, 0.9894336180360679, 1.019504458162858, 1.031824345179793, 1.02926596862519, 1.031231931498754, 0.9929783193494215, 1.022936654680617, 1.030727644359095),
    @(1.02, 1.022032777225756, 1.027149433735004, 0.9882828385668249, 1.017837680220131, 1.031628791838068, 1.0287207547942, 1.030762824586761, 0.9920501090198102, 1.021486749114238, 1.030181656261927),
    @(1.02, 1.021436676086645, 1.02672182460589, 0.9875604150241495, 1.016789608021511, 1.031504596826823, 1.028377289954112, 1.030467161601272, 0.9914670000341481, 1.020574531462177, 1.029837703662371),
    @(1.02, 1.021752678419797, 1.026948504520853, 0.9879432794643023, 1.01734522228031, 1.031570552702612, 1.028559430364201, 1.030623965980078, 0.991776070289318, 1.021058173964998, 1.030020102732734),
    @(1.02, 1.02299679624911, 1.02784100122195, 0.9894529299347244, 1.019532400534754, 1.031827602575562, 1.029275097885298, 1.031239783961406, 0.9929938892766442, 1.022960952597838, 1.0307367865838),
    @(1.02, 1.02444060819695, 1.028876826476721, 0.9912096547607049, 1.022070213526037, 1.032120476834377, 1.030102679556331, 1.031951273472153, 0.9944092447426414, 1.02516652614416, 1.03156554351573)
)

$nRows = 24
$left = New-Object "object[,]" $nRows,5
$right = New-Object "object[,]" $nRows,6
for ($i = 0; $i -lt $nRows; $i++) {
    $vals = $rowsData[$i]
    for ($j = 0; $j -lt 5; $j++) { $left[$i, $j] = $vals[$j] }
    for ($j = 0; $j -lt 6; $j++) { $right[$i, $j] = $vals[5 + $j] }
}

$ws.Range("B2:F25").Value = $left
$ws.Range("I2:N25").Value = $right

Write-Output "vm_pu.xlsx updated for 380 kV case"
